$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.873.35'
$ws.Range("E2").Value = '  +1.54%  '
$ws.Range("D3").Value = '2.602.35'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '519.74'
$ws.Range("E5").Value = '  +2.84%  '
$ws.Range("D6").Value = '154.79'
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '0.593'
$ws.Range("E8").Value = '  +3.00%  '
$ws.Range("D9").Value = '6.68'
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E10").Value = '  +2.97%  '
$ws.Range("D11").Value = '0.348'
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("E12").Value = '  +1.73%  '
$ws.Range("D13").Value = '3.057.59'
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").Value = '60.897.59'
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("D15").Value = '21.70'
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").Value = '2.607.08'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = '4.74'
$ws.Range("E18").Value = '  -1.69%  '
$ws.Range("D19").Value = '353.30'
$ws.Range("E19").Value = '  +2.52%  '
$ws.Range("D20").Value = '10.56'
$ws.Range("E20").Value = '  +1.99%  '
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").Value = '61.23'
$ws.Range("E23").Value = '  +2.10%  '
$ws.Range("D24").Value = '0.427'
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = '0.166'
$ws.Range("E25").Value = '  +0.58%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '2.718.15'
$ws.Range("E26").Value = '  +0.73%  '
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("D28").Value = '0.0₃0848'
$ws.Range("E28").Value = '  +0.93%  '
$ws.Range("D29").Value = '7.36'
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("D31").Value = '6.29'
$ws.Range("E31").Value = '  +9.43%  '
$ws.Range("D32").Value = '19.38'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("D34").Value = '149.23'
$ws.Range("E34").Value = '  -2.33%  '
$ws.Range("D35").Value = '4.22'
$ws.Range("E35").Value = '  +6.33%  '
$ws.Range("D36").Value = '0.931'
$ws.Range("E36").Value = '  +9.88%  '
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("D38").Value = '1.49'
$ws.Range("E38").Value = '  +2.07%  '
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '0.848'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Value = '3.79'
$ws.Range("E40").Value = '  +1.06%  '
$ws.Range("D41").Value = '36.45'
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("D42").Value = '286.72'
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '0.626'
$ws.Range("E43").Value = '  +2.05%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '0.101'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("D45").Value = '0.0560'
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = '19.57'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '4.87'
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0237'
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("D51").Value = '18.93'
$ws.Range("E51").Value = '  +7.69%  '
